# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) to the results table and refreshes a handful
# of prediction/error values that changed after the refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" column header (match the bold/centered/bordered look
#     already used by the other header cells) ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Label values for the first iteration block (rows 2-11) ---
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# --- Label values for the second iteration block (rows 12-21) ---
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# --- Refreshed prediction/error values (first iteration block) ---
$ws.Range("D3").Value = 0.3647985392704811
$ws.Range("E3").Value = 0.3647985392704811

$ws.Range("D4").Value = 0.3652367564456764
$ws.Range("E4").Value = 0.3652367564456764

$ws.Range("D7").Value = 0.5939353365730755
$ws.Range("E7").Value = 0.4060646634269245

$ws.Range("D8").Value = 0.6530953847192162
$ws.Range("E8").Value = 0.3469046152807838

$ws.Range("D11").Value = 0.7111764947679295
$ws.Range("E11").Value = 0.2888235052320705
